$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the full name text in A6 (add spaces) ---
$ws.Range("A6").Value = "David Gonzalo Cordon Fontecha"

# --- Re-point the hyperlinks on B6 / D6 so the "Cordondavid532@gmail.com"
#     mailto link is primarily associated with D6 (and B6 keeps a copy),
#     dropping the old custom display text in the process.
#     (Deleting any one hyperlink on the sheet clears the whole collection
#     in this runtime, so we rebuild all four from scratch, in the order
#     they should appear.)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:josemmp14@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:andresvillamizar_g@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:Cordondavid532@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:Cordondavid532@gmail.com")

# Restore the shared "Hipervinculo" style on the touched cells so no new
# style gets introduced by the hyperlink rebuild.
$ws.Range("D2").Style = "Hipervínculo"
$ws.Range("D3").Style = "Hipervínculo"
$ws.Range("D6").Style = "Hipervínculo"
$ws.Range("B6").Style = "Hipervínculo"

# --- Leave the cursor where it was when the workbook was last saved ---
$ws.Range("D7").Select()
